# Update ObjTables header rows: rename the "id=" attribute to "class="
# in cell A1 of every worksheet whose header uses the
# "!!ObjTables type='Data' id='...'" convention.
$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("A1")
    $text = $cell.Value()

    if ($text -ne $null -and $text.ToString().Contains("!!ObjTables type='Data' id='")) {
        $newText = $text.ToString().Replace("!!ObjTables type='Data' id='", "!!ObjTables type='Data' class='")
        $cell.Value = $newText
    }
}
